# Add a "{{ klasa_god }}" merge-field placeholder right after the
# "113-02/24-01/" text in the KLASA: paragraph, matching the same
# run-level formatting (Arial, 11pt/22 half-pt, hr-HR) already used
# there, and wrapping the inner "klasa_god" word in proofErr spell-check
# markers -- exactly like the existing "{{ dod }}" placeholder further
# down in the document.

$d = $word.ActiveDocument

$findRange = $d.Content
$findRange.Find.Execute("113-02/24-01/", $false, $false, $false, $false, `
                         $false, $true, 1, $false, "", 0)

$insertStart = $findRange.End
$insertRange = $d.Range($insertStart, $insertStart)

# A full "flat OPC" WordOpenXML package is needed so the new runs/proofErr
# markers keep their exact formatting (rFonts/sz/szCs/lang) instead of
# picking up Word's plain-text-insert defaults. We give the temporary
# paragraph the *same* paragraph identity/properties as the paragraph
# we're inserting into, because merging paragraphs (see below) keeps the
# formatting of the second paragraph's mark.
$xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="38DE985F" w14:textId="63191D6A" w:rsidR="003C7E33" w:rsidRDefault="003C7E33" w:rsidP="003C7E33">
            <w:pPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="hr-HR"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="hr-HR"/>
              </w:rPr>
              <w:t xml:space="preserve">{{ </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="hr-HR"/>
              </w:rPr>
              <w:t>klasa_god</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Arial" w:hAnsi="Arial"/>
                <w:sz w:val="22"/>
                <w:szCs w:val="22"/>
                <w:lang w:val="hr-HR"/>
              </w:rPr>
              <w:t xml:space="preserve"> }}</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xmlFragment)

# InsertXML inserted our runs as a brand-new paragraph (because the
# fragment contains a <w:p>), splitting the original "KLASA:" paragraph
# in two. Delete the paragraph mark it introduced to merge the new runs
# back into the original paragraph, so the text stays on one line just
# like the target diff shows.
$mark = $d.Range($insertStart, $insertStart + 1)
$mark.Delete()
